$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "oy, sandık, millet, seçim, parti"
$ws.Range("B2").Value = 26
$ws.Range("C2").Value = "millet, oy, buluş, sandık, seçim"

# Row 3
$ws.Range("A3").Value = "bayram, enstitü, kutlu, bereket, kahrol"
$ws.Range("B3").Value = 21
$ws.Range("C3").Value = "bayram, atatürk, türk, kutlu, mustafa"

# Row 4
$ws.Range("A4").Value = "cumhurbaşkan, aday, imza, ata, ittifak"
$ws.Range("B4").Value = 20
$ws.Range("C4").Value = "cumhurbaşkan, aday, ata, ittifak, dr"

# Row 5
$ws.Range("A5").Value = "tv, program, konuk, ekran, yayın"
$ws.Range("B5").Value = 17
$ws.Range("C5").Value = "yayın, program, canlı, konuk, sun"

# Row 6
$ws.Range("A6").Value = "atatürk, yaz, rahmet, yusuf, şehit"
$ws.Range("B6").Value = 13
$ws.Range("C6").Value = "şehit, rahmet, atatürk, an, dönüm"

# Row 7
$ws.Range("A7").Value = "türkiye, memleket, parti, yurttaş, oy"
$ws.Range("B7").Value = 11
$ws.Range("C7").Value = "türk, türkiye, milliyetçi, yüzyıl, millet"

# Row 8
$ws.Range("A8").Value = "sanatçı, iyi, memleket, fetö, sahte"
$ws.Range("B8").Value = 10
$ws.Range("C8").Value = "iyi, sanatçı, allah, iş, çık"

# Row 9
$ws.Range("A9").Value = "deprem, konut, yara, hatay, yurttaş"
$ws.Range("C9").Value = "konut, deprem, depremzede, temel, hastane"

# Row 10
$ws.Range("A10").Value = "teşekkür, yok, hazine, trabzon, samimiyet"
$ws.Range("C10").Value = "teşekkür, başkan, dernek, muhteşem, misafirperverlikleri"

# Row 11
$ws.Range("A11").Value = "cadde, sokak, esnaf, genç, meydan"
$ws.Range("B11").Value = 7
$ws.Range("C11").Value = "esnaf, cadde, ziyaret, genç, yoğun"

# Row 12
$ws.Range("A12").Value = "basın, açıkla, gazeteci, medya, temiz"
$ws.Range("C12").Value = "basın, açıkla, medya, uygula, cemiyet"

# Row 13
$ws.Range("C13").Value = "dadaş, öv, güzel, çocuk, erzurum"
